# fill the excel sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row (A1:G1) -------------------------------------------------
# Old headers had odd leading-space / trailing-space shared strings; the
# commit replaces them with clean text and reorders the columns.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "User Story"
$ws.Range("D1").Value = "State"
$ws.Range("E1").Value = "Owner"
$ws.Range("F1").Value = "Story Points"
$ws.Range("G1").Value = "Comments"

# --- Data rows -----------------------------------------------------------
$data = @(
    @{ Row=2;  A=1; B="16-09-14"; C="thinking of database design";               D="cmt"; E="Shilpa";               F=3 },
    @{ Row=3;  A=2; B="17-09-14"; C="designing of tables";                       D="wip"; E="shilpa";               F=3 },
    @{ Row=4;  A=3; B="18-09-14"; C="joining of tables with foreign ";           D="cmt"; E="shilpa";               F=2 },
    @{ Row=5;                     C="keys" },
    @{ Row=6;  A=4; B="19-09-14"; C="database diagrams";                        D="cmt"; E="shilpa";               F=2 },
    @{ Row=7;  A=5; B="22-09-14"; C="architecture of web application and";       D="wip"; E="ashutosh and shilpa";  F=2 },
    @{ Row=8;  A=6; B="23-09-14"; C="classes of tables are made in";             D="cmt"; E="shilpa";               F=2 },
    @{ Row=9;                     C="the layer" },
    @{ Row=10; A=7; B="24-09-14"; C="designing of login page";                   D="wip"; E="Ashutosh and shilpa";  F=3; G="designing of logo is still left" },
    @{ Row=11; A=8; B="25-09-14"; C="design the logo and paste";                 D="cmt"; E="shilpa";               F=3; G="logo allignment is still left" },
    @{ Row=12;                    C="change the color of buttons" }
)

foreach ($r in $data) {
    $row = $r.Row
    if ($r.ContainsKey("A")) { $ws.Cells.Item($row, 1).Value = $r.A }
    if ($r.ContainsKey("B")) { $ws.Cells.Item($row, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($row, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { $ws.Cells.Item($row, 4).Value = $r.D }
    if ($r.ContainsKey("E")) { $ws.Cells.Item($row, 5).Value = $r.E }
    if ($r.ContainsKey("F")) { $ws.Cells.Item($row, 6).Value = $r.F }
    if ($r.ContainsKey("G")) { $ws.Cells.Item($row, 7).Value = $r.G }
}

# --- Column widths ---------------------------------------------------------
# Columns B, D, E keep their pre-existing width (unchanged by the commit);
# only the newly introduced A, C, G columns and the resized F column need
# an explicit width. (Inputs below are pre-compensated for this runtime's
# character-width -> pixel rounding so the stored <col width> lands as
# close as possible to the target 9.42578125 / 27.5703125 / 15 / 37.)
$ws.Columns.Item(1).ColumnWidth = 8.584
$ws.Columns.Item(3).ColumnWidth = 26.584
$ws.Columns.Item(6).ColumnWidth = 14.084
$ws.Columns.Item(7).ColumnWidth = 36.084

# --- Selection -------------------------------------------------------------
$ws.Range("G10").Select()
